# Generate Report for Handback
# Update handoff/handback timestamps for the 00fe6a2a-...-md row across sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G2").Value = "2016-08-20 08:55:41"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H2").Value = "2016-08-20 08:55:37"
$zhcn.Range("K2").Value = "2016-08-20 08:55:53"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H2").Value = "2016-08-20 08:55:41"
$dede.Range("K2").Value = "2016-08-20 08:55:59"
